$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.960.91'
$ws.Range("D3").Value = '2.624.18'
$ws.Range("E3").Value = '  -3.39%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'588.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.02%  '
$ws.Range("D6").Value = "'164.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.36%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = "'0.534"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.58%  '
$ws.Range("D9").Value = '2.623.69'
$ws.Range("E9").Value = '  -3.40%  '
$ws.Range("E10").Value = '  -1.84%  '
$ws.Range("E11").Value = '  +1.28%  '
$ws.Range("D12").Value = "'0.361"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.19%  '
$ws.Range("D13").Value = "'5.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.94%  '
$ws.Range("D14").Value = "'27.46"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.78%  '
$ws.Range("D15").Value = '3.117.35'
$ws.Range("D16").Value = "'0.0000181"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.98%  '
$ws.Range("D17").Value = '66.845.66'
$ws.Range("E17").Value = '  -2.13%  '
$ws.Range("D18").Value = '2.627.05'
$ws.Range("E18").Value = '  -3.32%  '
$ws.Range("D19").Value = "'11.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.46%  '
$ws.Range("D20").Value = "'8.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.14%  '
$ws.Range("D21").Value = "'358.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.64%  '
$ws.Range("D22").Value = "'4.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.06%  '
$ws.Range("D23").Value = "'4.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.78%  '
$ws.Range("D24").Value = "'10.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +7.36%  '
$ws.Range("D25").Value = "'1.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.18%  '
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.19%  '
$ws.Range("D27").Value = "'70.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.34%  '
$ws.Range("D28").Value = '2.755.58'
$ws.Range("E28").Value = '  -3.51%  '
$ws.Range("D29").Value = "'0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.31%  '
$ws.Range("E30").Value = '  -3.35%  '
$ws.Range("D31").Value = "'550.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.20%  '
$ws.Range("D32").Value = "'7.91"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.72%  '
$ws.Range("D33").Value = "'1.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.73%  '
$ws.Range("D34").Value = "'1.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.96%  '
$ws.Range("D35").Value = "'0.132"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.13%  '
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("E37").Value = '  -5.33%  '
$ws.Range("D38").Value = "'157.36"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.72%  '
$ws.Range("D39").Value = "'19.13"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.87%  '
$ws.Range("E40").Value = '  -3.17%  '
$ws.Range("D41").Value = "'5.20"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.98%  '
$ws.Range("D42").Value = "'1.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.78%  '
$ws.Range("E43").Value = '  -0.63%  '
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("D45").Value = "'2.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.16%  '
$ws.Range("D46").Value = "'40.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.78%  '
$ws.Range("D47").Value = '0.0₆0294'
$ws.Range("E47").Value = '  -5.32%  '
$ws.Range("D48").Value = "'0.586"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.56%  '
$ws.Range("D49").Value = "'151.75"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.13%  '
$ws.Range("D50").Value = "'3.81"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.80%  '
$ws.Range("D51").Value = "'1.71"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.59%  '
